# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns with latest scrape values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to keep a purely textual value (no numeric auto-conversion),
    # then restore the default "Normal" style so no formatting residue is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "30.362.80"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "1.938.07"
$ws.Range("E3").Value = "  -2.97%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.21%  "
Set-TextValue $ws.Range("D5") "250.44"
$ws.Range("E5").Value = "  -1.65%  "
Set-TextValue $ws.Range("D6") "0.7243"
$ws.Range("E6").Value = "  -6.89%  "
$ws.Range("E7").Value = "  +0.17%  "
Set-TextValue $ws.Range("D8") "0.3333"
$ws.Range("E8").Value = "  -4.16%  "
$ws.Range("E9").Value = "  +2.12%  "
Set-TextValue $ws.Range("D10") "0.07230"
$ws.Range("E10").Value = "  +2.59%  "
Set-TextValue $ws.Range("D11") "0.8117"
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "1.937.52"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("E14").Value = "  -2.65%  "
Set-TextValue $ws.Range("D15") "94.30"
$ws.Range("E15").Value = "  -6.55%  "
Set-TextValue $ws.Range("D16") "14.99"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").Value = "30.360.85"
$ws.Range("E17").Value = "  -2.76%  "
Set-TextValue $ws.Range("D18") "0.000008255"
$ws.Range("E18").Value = "  +2.95%  "
Set-TextValue $ws.Range("D19") "248.92"
$ws.Range("E19").Value = "  -8.64%  "
Set-TextValue $ws.Range("D20") "5.910"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "2.189.77"
$ws.Range("E21").Value = "  -2.83%  "
Set-TextValue $ws.Range("D22") "1.001"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  +0.27%  "
Set-TextValue $ws.Range("D24") "6.954"
$ws.Range("E24").Value = "  -2.01%  "
Set-TextValue $ws.Range("D25") "9.754"
$ws.Range("E25").Value = "  -2.35%  "
Set-TextValue $ws.Range("D26") "163.31"
$ws.Range("E26").Value = "  -1.06%  "
Set-TextValue $ws.Range("D27") "2.394"
$ws.Range("E27").Value = "  +0.07%  "
Set-TextValue $ws.Range("D28") "19.28"
$ws.Range("E28").Value = "  -3.21%  "
Set-TextValue $ws.Range("D29") "0.1331"
$ws.Range("E29").Value = "  -7.75%  "
Set-TextValue $ws.Range("D30") "1.572"
$ws.Range("E30").Value = "  -1.52%  "
Set-TextValue $ws.Range("D31") "1.346"
$ws.Range("E31").Value = "  -1.59%  "
Set-TextValue $ws.Range("D32") "4.445"
$ws.Range("E32").Value = "  -3.30%  "
Set-TextValue $ws.Range("D33") "4.195"
$ws.Range("E33").Value = "  -5.19%  "
Set-TextValue $ws.Range("D34") "0.05203"
$ws.Range("E34").Value = "  -0.98%  "
Set-TextValue $ws.Range("D35") "1.290"
$ws.Range("E35").Value = "  +5.48%  "
Set-TextValue $ws.Range("D36") "0.7508"
$ws.Range("E36").Value = "  -4.49%  "
Set-TextValue $ws.Range("D37") "2.751"
$ws.Range("E37").Value = "  -0.56%  "
Set-TextValue $ws.Range("D38") "0.01980"
$ws.Range("E38").Value = "  -1.08%  "
Set-TextValue $ws.Range("D39") "2.842"
$ws.Range("E39").Value = "  -2.58%  "
Set-TextValue $ws.Range("D40") "80.80"
$ws.Range("E40").Value = "  +1.04%  "
Set-TextValue $ws.Range("D41") "6.456"
$ws.Range("E41").Value = "  -3.80%  "
Set-TextValue $ws.Range("D42") "0.4544"
$ws.Range("E42").Value = "  -2.91%  "
Set-TextValue $ws.Range("D43") "2.039"
$ws.Range("E43").Value = "  -2.99%  "
Set-TextValue $ws.Range("D44") "0.8486"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +0.17%  "
Set-TextValue $ws.Range("D46") "102.15"
$ws.Range("E46").Value = "  -2.37%  "
Set-TextValue $ws.Range("D47") "9.823"
$ws.Range("E47").Value = "  -1.77%  "
Set-TextValue $ws.Range("D48") "7.442"
$ws.Range("E48").Value = "  -3.05%  "
Set-TextValue $ws.Range("D49") "36.88"
$ws.Range("E49").Value = "  -0.75%  "
Set-TextValue $ws.Range("D50") "0.4194"
$ws.Range("E50").Value = "  -2.83%  "
Set-TextValue $ws.Range("D51") "2.869"
$ws.Range("E51").Value = "  +5.07%  "
